$wb = $excel.ActiveWorkbook

# "Ready for handoff" -> "In Translation" for every status cell across all
# three sheets (Overview has one status column per language: E & F;
# zh-cn / de-de each have a single Status column: C).
$wb.Worksheets.Item("Overview").Range("E2:F4").Value = "In Translation"
$wb.Worksheets.Item("zh-cn").Range("C2:C4").Value = "In Translation"
$wb.Worksheets.Item("de-de").Range("C2:C4").Value = "In Translation"

# The Status columns re-size (narrower, since "In Translation" renders
# narrower than "Ready for handoff") once the report is regenerated.
$wb.Worksheets.Item("Overview").Columns("E:F").ColumnWidth = 13.4101848602295
$wb.Worksheets.Item("zh-cn").Columns("C:C").ColumnWidth = 13.4101848602295
$wb.Worksheets.Item("de-de").Columns("C:C").ColumnWidth = 13.4101848602295
